# The uploaded workbook dropped the "IBAN" (bank account) column of fake
# account numbers that used to live in M2:M104 of Hoja1. The header in M1
# ("IBAN") is kept, only the per-row values are cleared out.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("M2:M104").ClearContents()
